$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Cells.Item(1,1).Value2 = "Datos actualizados a 25 de Agosto de 2020 a las 12:34"

# Update per-country numeric stats (B: Casos totales, C: Nuevos casos,
# D: Casos activos, E: Recuperados, F: Casos criticos, G: Muertes hoy, H: Muertes)

# Row 4: Estados Unidos
$ws.Cells.Item(4,2).Value2 = 5915911
$ws.Cells.Item(4,3).Value2 = 281
$ws.Cells.Item(4,4).Value2 = 3218514
$ws.Cells.Item(4,5).Value2 = 2516280
$ws.Cells.Item(4,7).Value2 = 3
$ws.Cells.Item(4,8).Value2 = 181117

# Row 14: Iran
$ws.Cells.Item(14,2).Value2 = 363363
$ws.Cells.Item(14,3).Value2 = 2213
$ws.Cells.Item(14,4).Value2 = 313058
$ws.Cells.Item(14,5).Value2 = 29404
$ws.Cells.Item(14,7).Value2 = 125
$ws.Cells.Item(14,8).Value2 = 20901

# Row 18: Banglades
$ws.Cells.Item(18,2).Value2 = 299628
$ws.Cells.Item(18,3).Value2 = 2545
$ws.Cells.Item(18,4).Value2 = 186756
$ws.Cells.Item(18,5).Value2 = 108844
$ws.Cells.Item(18,7).Value2 = 45
$ws.Cells.Item(18,8).Value2 = 4028

# Row 39: Oman
$ws.Cells.Item(39,2).Value2 = 84652
$ws.Cells.Item(39,3).Value2 = 143
$ws.Cells.Item(39,4).Value2 = 79147
$ws.Cells.Item(39,5).Value2 = 4863
$ws.Cells.Item(39,7).Value2 = 5
$ws.Cells.Item(39,8).Value2 = 642

# Row 42: Rumania
$ws.Cells.Item(42,2).Value2 = 80390
$ws.Cells.Item(42,3).Value2 = 1060
$ws.Cells.Item(42,4).Value2 = 35816
$ws.Cells.Item(42,5).Value2 = 41207
$ws.Cells.Item(42,7).Value2 = 58
$ws.Cells.Item(42,8).Value2 = 3367

# Row 86: Senegal
$ws.Cells.Item(86,2).Value2 = 13056
$ws.Cells.Item(86,3).Value2 = 43
$ws.Cells.Item(86,4).Value2 = 8715
$ws.Cells.Item(86,5).Value2 = 4067
$ws.Cells.Item(86,7).Value2 = 2
$ws.Cells.Item(86,8).Value2 = 274

# Row 91: Consejo Danes para los Refugiados
$ws.Cells.Item(91,2).Value2 = 9891
$ws.Cells.Item(91,3).Value2 = 49
$ws.Cells.Item(91,4).Value2 = 8972
$ws.Cells.Item(91,5).Value2 = 668

# Row 101: Finlandia
$ws.Cells.Item(101,2).Value2 = 7981
$ws.Cells.Item(101,3).Value2 = 43
$ws.Cells.Item(101,5).Value2 = 546

# Row 126: Sri Lanka
$ws.Cells.Item(126,4).Value2 = 2816
$ws.Cells.Item(126,5).Value2 = 131

# Row 175: Islas Feroe
$ws.Cells.Item(175,2).Value2 = 411
$ws.Cells.Item(175,3).Value2 = 1
$ws.Cells.Item(175,4).Value2 = 344
$ws.Cells.Item(175,5).Value2 = 67

# Row 185: Gibraltar
$ws.Cells.Item(185,2).Value2 = 256
$ws.Cells.Item(185,3).Value2 = 8
$ws.Cells.Item(185,5).Value2 = 53

# Rows 214/215: the source data reorders "Islas Malvinas" ahead of "Montserrat"
# (their shared-string entries swap position), so row 214 now carries Islas
# Malvinas' country name/figures and row 215 carries Montserrat's.
$ws.Cells.Item(214,1).Value2 = "Islas Malvinas"
$ws.Cells.Item(214,4).Value2 = 13
$ws.Cells.Item(214,8).Value2 = 0

$ws.Cells.Item(215,1).Value2 = "Montserrat"
$ws.Cells.Item(215,4).Value2 = 12
$ws.Cells.Item(215,8).Value2 = 1
